$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.462.43"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "3.810.71"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'669.60"
$ws.Range("E5").Value = "  +7.23%  "
$ws.Range("D6").Value = "'169.52"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").Value = "3.808.61"
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.525"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("E11").Value = "  +6.54%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "'35.71"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "4.451.40"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "3.808.50"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "70.477.54"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'11.81"
$ws.Range("E18").Value = "  +23.73%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'17.64"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "'476.69"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").Value = "'83.48"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("E25").Value = "  -3.85%  "
$ws.Range("D26").Value = "'12.21"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").Value = "'10.28"
$ws.Range("E27").Value = "  +2.62%  "
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "3.960.31"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("E31").Value = "  +7.45%  "
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("D33").Value = "'7.39"
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("D34").Value = "'29.54"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("E35").Value = "  +5.39%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("D38").Value = "3.766.07"
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "'3.37"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "'5.95"
$ws.Range("E41").Value = "  +2.45%  "
$ws.Range("D42").Value = "'0.962"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'2.11"
$ws.Range("E44").Value = "  +11.35%  "
$ws.Range("D46").Value = "'45.70"
$ws.Range("E46").Value = "  +6.06%  "
$ws.Range("D47").Value = "'158.99"
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("D48").Value = "'48.09"
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'1.42"
$ws.Range("E50").Value = "  +4.37%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000290"
$ws.Range("E51").Value = "  +5.85%  "
